# Bug Tracking System - minor fixes and DB update
# Adds two new defect rows (Jose's items) below the existing "Sprint initial
# story points..." row, assigns Jose to the existing row too, and keeps the
# trailing blank/total rows intact (shifted down).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert four new rows after row 7 (two data rows + two blank spacer rows),
# copying the formatting from row 7 so borders/number formats match the
# rest of the table.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(11).Insert()

$ws.Range("A7:F7").Copy()
$ws.Range("A8:F11").PasteSpecial(-4122)

# New defect row 9: "Remove Project list slider..." assigned to jose
$ws.Range("C9").Value() = "jose"

# Assign the existing open defect (row 7) and the new row 8 to Jose.
$ws.Range("C7").Value() = "Jose"
$ws.Range("C8").Value() = "Jose"

# New defect row 8: "Name Not changing..."
$ws.Range("B8").Value() = "Name Not changing on all files.  CompanyProperties business name should be on all headers not hard coded beter software"

# New defect row 9: "Remove Project list slider..."
$ws.Range("B9").Value() = "Remove Project list slider. Not used effectivly no longer makes sense."

$ws.Range("A8").Value() = 41604
$ws.Range("D8").Value() = 0.5
$ws.Rows.Item(8).RowHeight = 60

$ws.Range("A9").Value() = 41604
$ws.Range("D9").Value() = 0.5
$ws.Rows.Item(9).RowHeight = 30

# Update the selected cell to reflect where the user left off editing.
$ws.Range("A10").Select()
